$d = $word.ActiveDocument

# 1. Remove the LRI sentence from the summary paragraph text.
$lriSentence = "LRI (logisk ræsonnerings-indeks) blev målt til 67 (95% KI mellem 76-91), hvilket er langt under gennemsnittet. Denne score var 1. percentil, hvilket vil sige at 1% af børnene i norm-gruppen scorede lavere. "
$found1 = $d.Content.Find.Execute($lriSentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 2. Remove the LRI row from the results table.
$table = $d.Tables(1)
for ($i = 1; $i -le $table.Rows.Count; $i++) {
    $row = $table.Rows($i)
    $cellText = $row.Cells(1).Range.Text.TrimEnd([char]13, [char]7)
    if ($cellText -eq "LRI") {
        $row.Delete()
        break
    }
}

# 3. Remove the RSI ("RæsonneringsIndeks") recommendation block from the
#    "Anbefalinger" section, keeping the blank-line separator that precedes
#    the following "Arbejdshukommelses Indeks" block.
$br = [char]11
$rsiBlock = "" + $br + $br + "RæsonneringsIndeks måler evnen til at tænke logisk og abstrakt, udvikle nye tankemønstre, løse problemer på en kreativ måde, drage konklusioner og se mønstre og sammenhænge." + $br + "Indekset måler også overordnet visuel evne og evnen til at bearbejde informationer samtidigt." + $br + "Følgende støtte kan være relevant for en elev, som scorer lavt på RSI :" + $br + "- Lær eleven at bruge selvtale og verbal hukommelse til problemløsning" + $br + "- Lær barnet teknikke for problemløsning (fx at arbejde trin for trin, marker eller saml relevante dele til løsning af opgave)" + $br + "- Brug lister/procedurer til problemløsning" + $br + "- Overindlær gennem repetition" + $br + "- Undgå lange instruktioner og metaforer" + $br + "- Vær opmærksom på om der er vanskeligheder med at organisere og med sociale færdigheder"
$found3 = $d.Content.Find.Execute($rsiBlock, $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

Write-Output "sentence removed: $found1"
Write-Output "rsi block removed: $found3"
